# Add a new "Rectangle 367" text box (containing a "[Year]" placeholder)
# as a new direct child wps:wsp of the outer wpg:wgp group, right after the
# existing nested wpg:grpSp group closes - mirrors the upstream fix for
# "noop parent transformation" (zero chOff/chExt) handling on DOCX import.

$d = $word.ActiveDocument

# Locate the paragraph that holds the group drawing (wpg:wgp via a:graphic).
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -ne $null -and $para.Range.WordOpenXML -like "*wpg:wgp*") {
        $target = $para
        break
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs(1)
}

$range = $target.Range
$pkgXml = $range.WordOpenXML

# Pull the namespace declarations off the package's <w:document> root so the
# extracted <w:p> fragment remains well-formed once reinserted on its own.
$rootStart = $pkgXml.IndexOf("<w:document ")
$rootEnd = $pkgXml.IndexOf(">", $rootStart)
$rootTag = $pkgXml.Substring($rootStart, $rootEnd - $rootStart)
$nsParts = @()
foreach ($piece in $rootTag.Split(" ")) {
    if ($piece.StartsWith("xmlns:")) {
        $nsParts += $piece
    }
}
$nsAttr = [string]::Join(" ", $nsParts)

# Pull out the <w:p>...</w:p> fragment that contains the drawing.
$pStart = $pkgXml.IndexOf("<w:p ")
if ($pStart -lt 0) { $pStart = $pkgXml.IndexOf("<w:p>") }
$pEndMarker = "</w:p>"
$pEnd = $pkgXml.IndexOf($pEndMarker, $pStart) + $pEndMarker.Length
$pFrag = $pkgXml.Substring($pStart, $pEnd - $pStart)

# Re-declare the namespaces (a:, wps:, wpg:, ...) directly on <w:p> since it
# is being lifted out of the package and reinserted as a standalone fragment.
if ($pFrag.StartsWith("<w:p>")) {
    $pFrag = "<w:p " + $nsAttr + ">" + $pFrag.Substring(5)
} else {
    $pFrag = "<w:p " + $nsAttr + " " + $pFrag.Substring(5)
}

$newShape = @'
<wps:wsp><wps:cNvPr id="367" name="Rectangle 367"/><wps:cNvSpPr><a:spLocks noChangeArrowheads="1"/></wps:cNvSpPr><wps:spPr bwMode="auto"><a:xfrm><a:off x="7344" y="0"/><a:ext cx="4896" cy="3958"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/></wps:spPr><wps:txbx><w:txbxContent><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="96"/><w:szCs w:val="96"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="96"/><w:szCs w:val="96"/></w:rPr><w:t>[Year]</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" vert="horz" wrap="square" lIns="365760" tIns="182880" rIns="182880" bIns="182880" anchor="b" anchorCtr="0" upright="1"><a:noAutofit/></wps:bodyPr></wps:wsp>
'@

$marker = "</wpg:grpSp>"
$markerIdx = $pFrag.LastIndexOf($marker)
if ($markerIdx -lt 0) {
    throw "Could not find </wpg:grpSp> to anchor the new shape insertion"
}
$insertAt = $markerIdx + $marker.Length
$newPFrag = $pFrag.Substring(0, $insertAt) + $newShape + $pFrag.Substring($insertAt)

[void]$range.InsertXML($newPFrag)
